$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Website")

# Update the promo line text on the Website sheet.
$ws.Range("C7").Value = "You get `$20 for each new customer!"

# The "Cards that create new customers!" line is no longer bold.
$ws.Range("C6").Font.Bold = $false

# Set the page to portrait orientation (new pageSetup on this sheet).
$ws.PageSetup.Orientation = 1

# Move the cursor/selection to C15 as left by the author.
$null = $ws.Range("C15").Select()
